$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3000.2727
$ws.Range("I40").Value = 1857.7142
$ws.Range("J40").Value = 4999.75
$ws.Range("K40").Value = 1857.7142
$ws.Range("L40").Value = 4999.75
$ws.Range("M40").Value = -1682.7142
$ws.Range("N40").Value = -5349.75
$ws.Range("H132").Value = 3114.9375
$ws.Range("I132").Value = 2860.7556
$ws.Range("J132").Value = 6927.6665
$ws.Range("K132").Value = 8582.266799999999
$ws.Range("L132").Value = 20782.9995
$ws.Range("M132").Value = -6052.266799999999
$ws.Range("N132").Value = -25842.9995
$ws.Range("H137").Value = 249527.47
$ws.Range("I137").Value = 938150.25
$ws.Range("J137").Value = 5177.4517
$ws.Range("K137").Value = 2814450.75
$ws.Range("L137").Value = 15532.3551
$ws.Range("M137").Value = -2811900.75
$ws.Range("N137").Value = -20632.3551

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 2908.9565
$ws.Range("I74").Value = 2093.7334
$ws.Range("K74").Value = 2093.7334
$ws.Range("M74").Value = -1219.7334
$ws.Range("H77").Value = 2908.9565
$ws.Range("I77").Value = 2093.7334
$ws.Range("K77").Value = 10468.667
$ws.Range("M77").Value = -6100.667000000001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 30651.555
$ws.Range("I26").Value = 32457
$ws.Range("K26").Value = 32457
$ws.Range("M26").Value = -32165
$ws.Range("H94").Value = 12269.68
$ws.Range("I94").Value = 16778
$ws.Range("K94").Value = 16778
$ws.Range("M94").Value = -16327
$ws.Range("H99").Value = 26494.732
$ws.Range("I99").Value = 31701.75
$ws.Range("K99").Value = 31701.75
$ws.Range("M99").Value = -30203.75
$ws.Range("H134").Value = 9086.519
$ws.Range("I134").Value = 9486.916999999999
$ws.Range("K134").Value = 28460.751
$ws.Range("M134").Value = -25925.751

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H14").Value = 90006104
$ws.Range("I14").Value = 5877.5
$ws.Range("J14").Value = 150006260
$ws.Range("K14").Value = 5877.5
$ws.Range("L14").Value = 150006260
$ws.Range("M14").Value = -5707.5
$ws.Range("N14").Value = -150006600
$ws.Range("H31").Value = 3892
$ws.Range("I31").Value = 2665.6667
$ws.Range("J31").Value = 4028.2593
$ws.Range("K31").Value = 2665.6667
$ws.Range("L31").Value = 4028.2593
$ws.Range("M31").Value = -2370.6667
$ws.Range("N31").Value = -4618.2593
$ws.Range("H34").Value = 3892
$ws.Range("I34").Value = 2665.6667
$ws.Range("J34").Value = 4028.2593
$ws.Range("K34").Value = 2665.6667
$ws.Range("L34").Value = 4028.2593
$ws.Range("M34").Value = -2463.6667
$ws.Range("N34").Value = -4432.2593
$ws.Range("H132").Value = 27173.348
$ws.Range("I132").Value = 11711.632
$ws.Range("K132").Value = 35134.896
$ws.Range("M132").Value = -32604.896
$ws.Range("H134").Value = 3297315.5
$ws.Range("I134").Value = 4818130.5
$ws.Range("K134").Value = 14454391.5
$ws.Range("M134").Value = -14451856.5

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H11").Value = 125337.5
$ws.Range("J11").Value = 143235.72
$ws.Range("L11").Value = 429707.16
$ws.Range("N11").Value = -429987.16
$ws.Range("H18").Value = 3095.4285
$ws.Range("I18").Value = 3551.3333
$ws.Range("J18").Value = 2274.8
$ws.Range("K18").Value = 10653.9999
$ws.Range("L18").Value = 6824.400000000001
$ws.Range("M18").Value = -10484.9999
$ws.Range("N18").Value = -7162.400000000001
$ws.Range("H107").Value = 4076.4
$ws.Range("I107").Value = 948.3333
$ws.Range("K107").Value = 2844.9999
$ws.Range("M107").Value = -924.9998999999998
$ws.Range("H138").Value = 1002264.7
$ws.Range("I138").Value = 1113405.2
$ws.Range("K138").Value = 3340215.6
$ws.Range("M138").Value = -3335075.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 5624.75
$ws.Range("I102").Value = 6114
$ws.Range("K102").Value = 6114
$ws.Range("M102").Value = -4492
$ws.Range("H122").Value = 12277.421
$ws.Range("I122").Value = 14341.786
$ws.Range("J122").Value = 6497.2
$ws.Range("K122").Value = 43025.358
$ws.Range("L122").Value = 19491.6
$ws.Range("M122").Value = -40575.358
$ws.Range("N122").Value = -24391.6
$ws.Range("H132").Value = 4252.524
$ws.Range("I132").Value = 3359.7778
$ws.Range("K132").Value = 10079.3334
$ws.Range("M132").Value = -7549.3334
$ws.Range("H133").Value = 89990
$ws.Range("J133").Value = 89990
$ws.Range("L133").Value = 89990
$ws.Range("N133").Value = -100110

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2516.6843
$ws.Range("I46").Value = 1098.125
$ws.Range("J46").Value = 3548.3635
$ws.Range("K46").Value = 1098.125
$ws.Range("L46").Value = 3548.3635
$ws.Range("M46").Value = -910.125
$ws.Range("N46").Value = -3924.3635
$ws.Range("H81").Value = 0
$ws.Range("J81").Value = 0
$ws.Range("L81").Value = 0
$ws.Range("N81").ClearContents()
$ws.Range("H84").Value = 0
$ws.Range("J84").Value = 0
$ws.Range("L84").Value = 0
$ws.Range("N84").ClearContents()
$ws.Range("H100").Value = 5094.778
$ws.Range("I100").Value = 5231.6875
$ws.Range("K100").Value = 5231.6875
$ws.Range("M100").Value = -4690.6875
$ws.Range("H122").Value = 3334.606
$ws.Range("I122").Value = 2932.4827
$ws.Range("K122").Value = 8797.4481
$ws.Range("M122").Value = -6347.4481
$ws.Range("H127").Value = 142998960
$ws.Range("J127").Value = 165452.5
$ws.Range("L127").Value = 165452.5
$ws.Range("N127").Value = -175372.5
$ws.Range("H132").Value = 5462905
$ws.Range("I132").Value = 7801292.5
$ws.Range("K132").Value = 23403877.5
$ws.Range("M132").Value = -23401347.5
$ws.Range("H136").Value = 5592.6
$ws.Range("I136").Value = 4561.8
$ws.Range("K136").Value = 13685.4
$ws.Range("M136").Value = -11135.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 25478.9
$ws.Range("J81").Value = 4500
$ws.Range("L81").Value = 9000
$ws.Range("N81").Value = -11122
$ws.Range("H84").Value = 25478.9
$ws.Range("J84").Value = 4500
$ws.Range("L84").Value = 45000
$ws.Range("N84").Value = -55608
$ws.Range("H107").Value = 12965
$ws.Range("I107").Value = 1786.3636
$ws.Range("K107").Value = 5359.0908
$ws.Range("M107").Value = -3439.0908
$ws.Range("H122").Value = 3908.6904
$ws.Range("I122").Value = 2544.2
$ws.Range("K122").Value = 7632.599999999999
$ws.Range("M122").Value = -5182.599999999999
$ws.Range("H132").Value = 8596.758
$ws.Range("I132").Value = 9761.405000000001
$ws.Range("J132").Value = 6873.08
$ws.Range("K132").Value = 29284.215
$ws.Range("L132").Value = 20619.24
$ws.Range("M132").Value = -26754.215
$ws.Range("N132").Value = -25679.24
